$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Component Data")
$lo = $ws.ListObjects.Item("ComponentData")
$col = $lo.ListColumns.Item("Efficiency")
$col.Name = "Transmission Loss"
